# Weekly NYPD CompStat refresh: new commissioner name, new volume/week
# numbers, new reporting dates, and a fresh week of crime-complaint figures
# (the "0" placeholder Murder-week-to-date string cell becomes a real
# number, and the table gains the "Transit" precinct row that was
# previously missing numbers).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Masthead -------------------------------------------------------
$ws.Range("M6").Value = "Edward A. Caban"
$ws.Range("A8").Value = "Volume 30   Number  27"
$ws.Range("C9").Value = "Report Covering the Week  7/3/2023  Through  7/9/2023"

# --- Weekly crime-complaint table (rows 14-30, columns C:N) ---------
# C = current week, D = year-ago week, E = %chg
# F = 28-day, G = 28-day year-ago, H = %chg
# I = YTD, J = YTD year-ago, K = %chg
# L/M/N = 2-year / 13-year / 30-year %chg
$weeklyData = @(
    @{ Row = 14; Vals = @(1, 1, 0, 1, 10, -90, 11, 25, -56, -54.166666666666, -67.647058823529, -84.722222222222) }
    @{ Row = 15; Vals = @(1, 3, -66.666666666666, 8, 19, -57.894736842105, 73, 105, -30.476190476190, -17.045454545454, 2.816901408450, -62.755102040816) }
    @{ Row = 16; Vals = @(28, 34, -17.647058823529, 99, 116, -14.655172413793, 662, 726, -8.815426997245, 21.915285451197, -39.042357274401, -82.960102960103) }
    @{ Row = 17; Vals = @(72, 84, -14.285714285714, 259, 271, -4.428044280442, 1435, 1412, 1.628895184135, 19.384359400998, 67.444574095682, -22.558014031300) }
    @{ Row = 18; Vals = @(20, 33, -39.393939393939, 79, 103, -23.300970873786, 568, 563, 0.888099467140, 22.413793103448, -46.313799621928, -87.976291278577) }
    @{ Row = 19; Vals = @(74, 79, -6.329113924050, 275, 327, -15.902140672782, 1790, 1942, -7.826982492276, 47.082990961380, 26.056338028169, -57.033125300048) }
    @{ Row = 20; Vals = @(38, 29, 31.034482758620, 158, 105, 50.476190476190, 899, 808, 11.262376237623, 80.522088353413, 6.516587677725, -90.896202531645) }
    @{ Row = 21; Vals = @(234, 263, -11.026615969581, 879, 951, -7.570977917981, 5438, 5581, -2.562264827091, 34.737363726461, 1.266294227188, -78.046909692785) }
    @{ Row = 22; Vals = @(6, 1, 500, 12, 5, 140, 62, 52, 19.230769230769, 31.914893617021, 0, "***.*") }
    @{ Row = 23; Vals = @(1, 4, -75, 17, 15, 13.333333333333, 120, 115, 4.347826086956, 10.091743119266, 48.148148148148, "***.*") }
    @{ Row = 24; Vals = @(167, 238, -29.831932773109, 740, 873, -15.234822451317, 4805, 5054, -4.926790660862, 43.518518518518, 52.153261557948, "***.*") }
    @{ Row = 25; Vals = @(95, 114, -16.666666666666, 379, 445, -14.831460674157, 2325, 2075, 12.048192771084, 36.124121779859, -3.044203502919, "***.*") }
    @{ Row = 26; Vals = @(1, 4, -75, 12, 29, -58.620689655172, 134, 170, -21.176470588235, -4.285714285714, "***.*", "***.*") }
    @{ Row = 27; Vals = @(8, 5, 60, 29, 33, -12.121212121212, 221, 223, -0.896860986547, 9.405940594059, "***.*", "***.*") }
    @{ Row = 28; Vals = @(4, 14, -71.428571428571, 11, 32, -65.625, 56, 97, -42.268041237113, -44.554455445544, -48.148148148148, -77.6) }
    @{ Row = 29; Vals = @(2, 8, -75, 7, 18, -61.111111111111, 40, 74, -45.945945945945, -53.488372093023, -53.488372093023, -82.683982683982) }
    @{ Row = 30; Vals = @(2, 1, 100, 6, 4, 50, 33, 19, 73.684210526315, 153.846153846154, "***.*", "***.*") }
)

foreach ($item in $weeklyData) {
    $r = $item.Row
    $col = 3
    foreach ($v in $item.Vals) {
        $ws.Cells.Item($r, $col).Value = $v
        $col = $col + 1
    }
}

# C14 ("Murder" week-to-date) used to hold the literal text placeholder
# "0"; it now carries a real number, so pick up the numeric formatting
# used by the rest of the column.
$ws.Range("C14").NumberFormat = $ws.Range("C15").NumberFormat
